$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-11-30 Saturday" "2024-12-01 Sunday"

Replace-Text "36×76=" "20×75="
Replace-Text "93×92=" "13×70="
Replace-Text "93×31=" "27×77="
Replace-Text "20×23=" "27×40="
Replace-Text "28×53=" "82×29="
Replace-Text "23×98=" "54×23="
Replace-Text "30×49=" "18×64="
Replace-Text "66×51=" "19×23="
Replace-Text "79×82=" "22×17="
Replace-Text "54×76=" "12×68="
Replace-Text "48×93=" "16×55="
Replace-Text "54×39=" "69×60="
Replace-Text "65×44=" "77×65="
Replace-Text "26×18=" "92×16="
Replace-Text "24×35=" "97×42="
Replace-Text "45×35=" "37×59="
Replace-Text "12×86=" "65×99="
Replace-Text "93×28=" "71×75="
Replace-Text "18×77=" "24×41="
Replace-Text "51×88=" "75×85="
Replace-Text "77×85=" "37×19="
Replace-Text "33×78=" "56×70="
Replace-Text "79×41=" "61×81="
Replace-Text "30×69=" "57×24="
Replace-Text "89×56=" "47×48="
